$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1563.6666
$ws.Range("I19").Value = 890.5
$ws.Range("J19").Value = 1900.25
$ws.Range("K19").Value = 890.5
$ws.Range("L19").Value = 1900.25
$ws.Range("M19").Value = -715.5
$ws.Range("N19").Value = -2250.25

$ws.Range("H74").Value = 3924.875
$ws.Range("I74").Value = 3857
$ws.Range("K74").Value = 3857
$ws.Range("M74").Value = -2921

$ws.Range("H77").Value = 3924.875
$ws.Range("I77").Value = 3857
$ws.Range("K77").Value = 19285
$ws.Range("M77").Value = -14605

$ws.Range("H100").Value = 9260516
$ws.Range("I100").Value = 12346240
$ws.Range("J100").Value = 3342.7778
$ws.Range("K100").Value = 12346240
$ws.Range("L100").Value = 3342.7778
$ws.Range("M100").Value = -12345699
$ws.Range("N100").Value = -4424.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2970.25
$ws.Range("I45").Value = 2828.818
$ws.Range("K45").Value = 2828.818
$ws.Range("M45").Value = -2451.818

$ws.Range("H61").Value = 1308085.2
$ws.Range("I61").Value = 1667460.2
$ws.Range("K61").Value = 1667460.2
$ws.Range("M61").Value = -1667248.2

$ws.Range("H74").Value = 619.6429000000001
$ws.Range("I74").Value = 574.1667
$ws.Range("J74").Value = 892.5
$ws.Range("K74").Value = 574.1667
$ws.Range("L74").Value = 892.5
$ws.Range("M74").Value = 299.8333
$ws.Range("N74").Value = -2640.5

$ws.Range("H77").Value = 619.6429000000001
$ws.Range("I77").Value = 574.1667
$ws.Range("J77").Value = 892.5
$ws.Range("K77").Value = 2870.8335
$ws.Range("L77").Value = 4462.5
$ws.Range("M77").Value = 1497.1665
$ws.Range("N77").Value = -13198.5

$ws.Range("H97").Value = 1554.8
$ws.Range("I97").Value = 964.2857
$ws.Range("K97").Value = 964.2857
$ws.Range("M97").Value = -468.2857

$ws.Range("H122").Value = 58121
$ws.Range("I122").Value = 2841.5
$ws.Range("J122").Value = 500357
$ws.Range("K122").Value = 8524.5
$ws.Range("L122").Value = 1501071
$ws.Range("M122").Value = -6074.5
$ws.Range("N122").Value = -1505971

$ws.Range("H132").Value = 26782.486
$ws.Range("I132").Value = 688.86206
$ws.Range("K132").Value = 2066.58618
$ws.Range("M132").Value = 463.4138199999998

$ws.Range("H136").Value = 1308085.2
$ws.Range("I136").Value = 1667460.2
$ws.Range("K136").Value = 5002380.6
$ws.Range("M136").Value = -4999830.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1081.0667
$ws.Range("I94").Value = 1212.5714
$ws.Range("J94").Value = 966
$ws.Range("K94").Value = 1212.5714
$ws.Range("L94").Value = 966
$ws.Range("M94").Value = -761.5714
$ws.Range("N94").Value = -1868

$ws.Range("H99").Value = 1211.4546
$ws.Range("I99").Value = 850.64703
$ws.Range("J99").Value = 2438.2
$ws.Range("K99").Value = 850.64703
$ws.Range("L99").Value = 2438.2
$ws.Range("M99").Value = 647.35297
$ws.Range("N99").Value = -5434.2

$ws.Range("H107").Value = 1392.2222
$ws.Range("I107").Value = 1014.3333
$ws.Range("J107").Value = 2148
$ws.Range("K107").Value = 1014.3333
$ws.Range("L107").Value = 2148
$ws.Range("M107").Value = 905.6667
$ws.Range("N107").Value = -5988

$ws.Range("H134").Value = 6665.147
$ws.Range("I134").Value = 2252.1538
$ws.Range("K134").Value = 6756.4614
$ws.Range("M134").Value = -4221.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4809433
$ws.Range("I31").Value = 6411626.5
$ws.Range("J31").Value = 2853.7693
$ws.Range("K31").Value = 6411626.5
$ws.Range("L31").Value = 2853.7693
$ws.Range("M31").Value = -6411331.5
$ws.Range("N31").Value = -3443.7693

$ws.Range("H34").Value = 4809433
$ws.Range("I34").Value = 6411626.5
$ws.Range("J34").Value = 2853.7693
$ws.Range("K34").Value = 6411626.5
$ws.Range("L34").Value = 2853.7693
$ws.Range("M34").Value = -6411424.5
$ws.Range("N34").Value = -3257.7693

$ws.Range("H58").Value = 1457.2727
$ws.Range("I58").Value = 1405
$ws.Range("J58").Value = 1980
$ws.Range("K58").Value = 1405
$ws.Range("L58").Value = 1980
$ws.Range("M58").Value = -1202
$ws.Range("N58").Value = -2386

$ws.Range("H86").Value = 47602
$ws.Range("I86").Value = 15758.143
$ws.Range("J86").Value = 64748.69
$ws.Range("K86").Value = 15758.143
$ws.Range("L86").Value = 64748.69
$ws.Range("M86").Value = -14635.143
$ws.Range("N86").Value = -66994.69

$ws.Range("H89").Value = 47602
$ws.Range("I89").Value = 15758.143
$ws.Range("J89").Value = 64748.69
$ws.Range("K89").Value = 78790.715
$ws.Range("L89").Value = 323743.45
$ws.Range("M89").Value = -73174.715
$ws.Range("N89").Value = -334975.45

$ws.Range("H99").Value = 1519.9578
$ws.Range("I99").Value = 1497.2273
$ws.Range("J99").Value = 1820
$ws.Range("K99").Value = 1497.2273
$ws.Range("L99").Value = 1820
$ws.Range("M99").Value = 0.7726999999999862
$ws.Range("N99").Value = -4816

$ws.Range("H126").Value = 1519.9578
$ws.Range("I126").Value = 1497.2273
$ws.Range("J126").Value = 1820
$ws.Range("K126").Value = 4491.6819
$ws.Range("L126").Value = 5460
$ws.Range("M126").Value = -2021.6819
$ws.Range("N126").Value = -10400

$ws.Range("H132").Value = 39119.11
$ws.Range("I132").Value = 1307.6666
$ws.Range("K132").Value = 3922.9998
$ws.Range("M132").Value = -1392.9998

$ws.Range("H134").Value = 1720.1818
$ws.Range("I134").Value = 1148.9375
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 3446.8125
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -911.8125
$ws.Range("N134").Value = -65070

$ws.Range("H136").Value = 1457.2727
$ws.Range("I136").Value = 1405
$ws.Range("J136").Value = 1980
$ws.Range("K136").Value = 4215
$ws.Range("L136").Value = 5940
$ws.Range("M136").Value = -1665
$ws.Range("N136").Value = -11040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1247.1818
$ws.Range("I97").Value = 980.0909
$ws.Range("J97").Value = 1781.3636
$ws.Range("K97").Value = 980.0909
$ws.Range("L97").Value = 1781.3636
$ws.Range("M97").Value = -484.0909
$ws.Range("N97").Value = -2773.3636

$ws.Range("H102").Value = 1871.5264
$ws.Range("I102").Value = 1450.6666
$ws.Range("J102").Value = 3449.75
$ws.Range("K102").Value = 1450.6666
$ws.Range("L102").Value = 3449.75
$ws.Range("M102").Value = 171.3334
$ws.Range("N102").Value = -6693.75

$ws.Range("H113").Value = 1383.25
$ws.Range("I113").Value = 899.125
$ws.Range("J113").Value = 2351.5
$ws.Range("K113").Value = 899.125
$ws.Range("L113").Value = 2351.5
$ws.Range("M113").Value = 1270.875
$ws.Range("N113").Value = -6691.5

$ws.Range("H122").Value = 1620.8096
$ws.Range("I122").Value = 1548.0667
$ws.Range("K122").Value = 4644.2001
$ws.Range("M122").Value = -2194.2001

$ws.Range("H126").Value = 8334805.5
$ws.Range("I126").Value = 1386.5834
$ws.Range("J126").Value = 20834934
$ws.Range("K126").Value = 4159.7502
$ws.Range("L126").Value = 62504802
$ws.Range("M126").Value = -1689.7502
$ws.Range("N126").Value = -62509742

$ws.Range("H132").Value = 368359.78
$ws.Range("I132").Value = 44755.695
$ws.Range("J132").Value = 1431630.2
$ws.Range("K132").Value = 134267.085
$ws.Range("L132").Value = 4294890.6
$ws.Range("M132").Value = -131737.085
$ws.Range("N132").Value = -4299950.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2527.75
$ws.Range("I7").Value = 1533.3
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 1533.3
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -1421.3
$ws.Range("N7").Value = -7724

$ws.Range("H40").Value = 2003.0588
$ws.Range("I40").Value = 1894.3334
$ws.Range("J40").Value = 2125.375
$ws.Range("K40").Value = 1894.3334
$ws.Range("L40").Value = 2125.375
$ws.Range("M40").Value = -1758.3334
$ws.Range("N40").Value = -2397.375

$ws.Range("H61").Value = 1580.6786
$ws.Range("I61").Value = 1271.7894
$ws.Range("J61").Value = 2232.7778
$ws.Range("K61").Value = 1271.7894
$ws.Range("L61").Value = 2232.7778
$ws.Range("M61").Value = -1069.7894
$ws.Range("N61").Value = -2636.7778

$ws.Range("H93").Value = 1016.8333
$ws.Range("I93").Value = 944.1875
$ws.Range("J93").Value = 1162.125
$ws.Range("K93").Value = 944.1875
$ws.Range("L93").Value = 1162.125
$ws.Range("M93").Value = 303.8125
$ws.Range("N93").Value = -3658.125

$ws.Range("H100").Value = 2046.2222
$ws.Range("I100").Value = 2014.9333
$ws.Range("J100").Value = 2202.6667
$ws.Range("K100").Value = 2014.9333
$ws.Range("L100").Value = 2202.6667
$ws.Range("M100").Value = -1473.9333
$ws.Range("N100").Value = -3284.6667

$ws.Range("H113").Value = 1580.6786
$ws.Range("I113").Value = 1271.7894
$ws.Range("J113").Value = 2232.7778
$ws.Range("K113").Value = 1271.7894
$ws.Range("L113").Value = 2232.7778
$ws.Range("M113").Value = 898.2106000000001
$ws.Range("N113").Value = -6572.7778

$ws.Range("H122").Value = 1779.1111
$ws.Range("I122").Value = 1752
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 5256
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -2806
$ws.Range("N122").Value = -10400.0002

$ws.Range("H126").Value = 2527.75
$ws.Range("I126").Value = 1533.3
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 4599.9
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -2129.9
$ws.Range("N126").Value = -27440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1708.68
$ws.Range("I122").Value = 1590.8667
$ws.Range("J122").Value = 1885.4
$ws.Range("K122").Value = 4772.6001
$ws.Range("L122").Value = 5656.200000000001
$ws.Range("M122").Value = -2322.6001
$ws.Range("N122").Value = -10556.2

$ws.Range("H132").Value = 76413144
$ws.Range("I132").Value = 113000860
$ws.Range("J132").Value = 3237706.2
$ws.Range("K132").Value = 339002580
$ws.Range("L132").Value = 9713118.600000001
$ws.Range("M132").Value = -339000050
$ws.Range("N132").Value = -9718178.600000001

$ws.Range("H136").Value = 44692.957
$ws.Range("I136").Value = 63520.812
$ws.Range("J136").Value = 1657.8572
$ws.Range("K136").Value = 190562.436
$ws.Range("L136").Value = 4973.571599999999
$ws.Range("M136").Value = -188012.436
$ws.Range("N136").Value = -10073.5716
